# Slide 1 ("History"): the first bullet line currently reads
#   "Who ? \tLinux Torvalds (Invented Linux)"
# as a single run. The author split it into three runs, fixing the
# typo "Linux" -> "Linus" (Linus Torvalds) along the way:
#   "Who ? \t" + "Linus " + "Torvalds (Invented Linux)"
#
# Re-create that by editing only the substring that needs to change
# ("Linux " -> "Linus ") via TextRange.Characters(start, length), which
# naturally splits the run around the edited span, leaving the
# untouched "Who ? \t" and "Torvalds (Invented Linux)" portions intact.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item("Content Placeholder 4")
$tf = $shape.TextFrame
$tr = $tf.TextRange

$firstPara = $tr.Paragraphs(1, 1)

# "Who ? \t" is 7 characters, so the word to fix ("Linux ") starts at
# character 8 and is 6 characters long (includes the trailing space).
$target = $firstPara.Characters(8, 6)
$target.Text = "Linus "
